# Applies the "Added minimum operating point to spine and jupyter" edit.
#
# Summary of changes:
#  1. Sheet "Definition": the node rows (A7:A18) are reordered (their
#     Category in column B stays "node" throughout, only the names move).
#  2. Sheet "Nodes": rows 2-13 (all columns A-G) are reordered to the same
#     node order as "Definition" - the full per-node data (balance_type,
#     has_state, node_state_cap, frac_state_loss, node_slack_penalty)
#     travels together with its node name.
#  3. Sheet "Object__to_from_node": the unit__to_node relationship between
#     Solar_Plant_Kasso and Power_Kasso gets a new parameter
#     "minimum_operating_point" = 0.2, plus new "ramp_up_limit" (0.3) and
#     "ramp_down_limit" (0.1) parameters, and its "shut_down_limit" value
#     changes from 0.8 to 0.2. The rest of the relationship rows shift
#     down to make room (sheet grows from 22 to 25 data rows).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Definition" sheet - reorder node names in A7:A18 (column B is
#    always "node" for this block, so it does not need to change).
# ---------------------------------------------------------------------
$wsDef = $wb.Worksheets.Item("Definition")

$defNodeOrder = @(
    "District_Heating",
    "Power_Kasso",
    "Raw_Methanol",
    "Water",
    "E-Methanol_Kasso",
    "Hydrogen_Kasso",
    "Waste_Heat",
    "Hydrogen_storage_Kasso",
    "Vaporized_Carbon_Dioxide",
    "Power_Wholesale",
    "Carbon_Dioxide",
    "E-Methanol_storage_Kasso"
)

$r = 7
foreach ($name in $defNodeOrder) {
    $wsDef.Cells.Item($r, 1).Value = $name
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) "Nodes" sheet - reorder rows 2-13 (columns A-G) so each node's full
#    set of parameters travels with it, in the same order as "Definition".
#    Columns: A=Object_Name, B=Category, C=balance_type, D=has_state,
#             E=node_state_cap, F=frac_state_loss, G=node_slack_penalty
# ---------------------------------------------------------------------
$wsNodes = $wb.Worksheets.Item("Nodes")

# Each entry: Name, Category, balance_type, has_state, node_state_cap, frac_state_loss, node_slack_penalty
# $null means the cell should be blank.
$nodesData = @(
    @("District_Heating",          "node", "balance_type_none", $null,  $null,   $null, $null),
    @("Power_Kasso",                "node", "balance_type_node", $null,  $null,   $null, 100000),
    @("Raw_Methanol",                "node", "balance_type_node", $null,  $null,   $null, 100000),
    @("Water",                       "node", "balance_type_none", $null,  $null,   $null, $null),
    @("E-Methanol_Kasso",            "node", "balance_type_node", $null,  $null,   $null, 100000),
    @("Hydrogen_Kasso",               "node", "balance_type_node", $null,  $null,   $null, 100000),
    @("Waste_Heat",                   "node", "balance_type_node", $null,  $null,   $null, $null),
    @("Hydrogen_storage_Kasso",       "node", "balance_type_node", "true", 100000,  0,     100000),
    @("Vaporized_Carbon_Dioxide",     "node", "balance_type_node", $null,  $null,   $null, 100000),
    @("Power_Wholesale",              "node", "balance_type_none", $null,  $null,   $null, $null),
    @("Carbon_Dioxide",               "node", "balance_type_none", $null,  $null,   $null, $null),
    @("E-Methanol_storage_Kasso",     "node", "balance_type_node", "true", 100000,  0,     100000)
)

$r = 2
foreach ($row in $nodesData) {
    $wsNodes.Cells.Item($r, 1).Value = $row[0]
    $wsNodes.Cells.Item($r, 2).Value = $row[1]
    $wsNodes.Cells.Item($r, 3).Value = $row[2]

    if ($null -eq $row[3]) { $wsNodes.Cells.Item($r, 4).Value = "" } else { $wsNodes.Cells.Item($r, 4).Value = $row[3] }
    if ($null -eq $row[4]) { $wsNodes.Cells.Item($r, 5).Value = "" } else { $wsNodes.Cells.Item($r, 5).Value = $row[4] }
    if ($null -eq $row[5]) { $wsNodes.Cells.Item($r, 6).Value = "" } else { $wsNodes.Cells.Item($r, 6).Value = $row[5] }
    if ($null -eq $row[6]) { $wsNodes.Cells.Item($r, 7).Value = "" } else { $wsNodes.Cells.Item($r, 7).Value = $row[6] }

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3) "Object__to_from_node" sheet - rebuild rows 2-25.
#    Columns: A=relationship_class_name, B=object_class, C=object_name,
#             D=node, E=parameter_name, F=value
# ---------------------------------------------------------------------
$wsRel = $wb.Worksheets.Item("Object__to_from_node")

$relData = @(
    @("unit__to_node",        "unit",       "Solar_Plant_Kasso",          "Power_Kasso",               "unit_capacity",            304),
    @("unit__to_node",        "unit",       "Solar_Plant_Kasso",          "Power_Kasso",               "ramp_up_limit",            0.3),
    @("unit__to_node",        "unit",       "Solar_Plant_Kasso",          "Power_Kasso",               "ramp_down_limit",          0.1),
    @("unit__to_node",        "unit",       "Solar_Plant_Kasso",          "Power_Kasso",               "start_up_limit",           0.5),
    @("unit__to_node",        "unit",       "Solar_Plant_Kasso",          "Power_Kasso",               "shut_down_limit",          0.2),
    @("unit__to_node",        "unit",       "Solar_Plant_Kasso",          "Power_Kasso",               "minimum_operating_point",  0.2),
    @("unit__from_node",      "unit",       "Electrolyzer",                "Power_Kasso",               "unit_capacity",            52),
    @("unit__from_node",      "unit",       "Electrolyzer",                "Power_Kasso",               "vom_cost",                 1),
    @("unit__to_node",        "unit",       "CO2_Vaporizer",               "Vaporized_Carbon_Dioxide",  "unit_capacity",            100),
    @("unit__to_node",        "unit",       "Destilation_Tower",           "E-Methanol_Kasso",          "unit_capacity",            52),
    @("unit__to_node",        "unit",       "Methanol_Reactor",            "Raw_Methanol",               "unit_capacity",            100),
    @("unit__to_node",        "unit",       "Methanol_Reactor",            "Raw_Methanol",               "ramp_up_limit",            0.5),
    @("unit__to_node",        "unit",       "Methanol_Reactor",            "Raw_Methanol",               "ramp_down_limit",          0.5),
    @("unit__to_node",        "unit",       "Methanol_Reactor",            "Waste_Heat",                 "unit_capacity",            100),
    @("connection__from_node","connection", "power_line_Wholesale_Kasso",  "Power_Wholesale",            "connection_capacity",      1000),
    @("connection__to_node",  "connection", "power_line_Wholesale_Kasso",  "Power_Kasso",                "connection_capacity",      1000),
    @("connection__from_node","connection", "power_line_Wholesale_Kasso",  "Power_Kasso",                "connection_capacity",      1000),
    @("connection__to_node",  "connection", "power_line_Wholesale_Kasso",  "Power_Wholesale",            "connection_capacity",      1000),
    @("connection__to_node",  "connection", "pipeline_storage_hydrogen",   "Hydrogen_storage_Kasso",     "connection_capacity",      1000),
    @("connection__from_node","connection", "pipeline_storage_hydrogen",   "Hydrogen_storage_Kasso",     "connection_capacity",      1000),
    @("connection__to_node",  "connection", "pipeline_storage_e-methanol", "E-Methanol_storage_Kasso",   "connection_capacity",      1000),
    @("connection__from_node","connection", "pipeline_storage_e-methanol", "E-Methanol_storage_Kasso",   "connection_capacity",      1000),
    @("connection__from_node","connection", "pipeline_District_Heating",   "Waste_Heat",                 "connection_capacity",      1000),
    @("connection__to_node",  "connection", "pipeline_District_Heating",   "District_Heating",           "connection_capacity",      1000)
)

$r = 2
foreach ($row in $relData) {
    $wsRel.Cells.Item($r, 1).Value = $row[0]
    $wsRel.Cells.Item($r, 2).Value = $row[1]
    $wsRel.Cells.Item($r, 3).Value = $row[2]
    $wsRel.Cells.Item($r, 4).Value = $row[3]
    $wsRel.Cells.Item($r, 5).Value = $row[4]
    $wsRel.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

Write-Host "Edit applied successfully"
